# Add withdrawal report to get correct buy-dates of forex for potential
# transactions in the future.
#
# Re-derives the "Foreign Currencies" forex gain/loss rows against the
# withdrawal-based buy dates instead of the original purchase dates, and
# drops the rows that are no longer part of the recomputed withdrawal
# matching. The "ELSTER - Summary" sheet's forex total is refreshed to
# match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foreign Currencies")

# --- Row 6's old "Sell Date" (s23) becomes row 5's new "Buy Date";     ---
# --- do this copy before row 5 itself is overwritten below.           ---
$ws.Cells.Item(5, 3).Copy()
$ws.Cells.Item(6, 3).PasteSpecial(-4163)

# --- Row 4's "Buy Date" (s22) also becomes row 5's new "Buy Date".     ---
$ws.Cells.Item(4, 3).Copy()
$ws.Cells.Item(5, 3).PasteSpecial(-4163)

# --- Rows 5 & 6 pick up the later "2022-12-01" (s45) withdrawal date   ---
# --- already used as the "Sell Date" for rows 7 & 8.                  ---
$ws.Cells.Item(7, 4).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4163)
$ws.Cells.Item(8, 4).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)

$excel.CutCopyMode = 0

# --- Update the recomputed quantities / rates / gains for rows 4-6 ------
$ws.Cells.Item(4, 2).Value = 2582.03
$ws.Cells.Item(4, 7).Value = 57.39

$ws.Cells.Item(5, 2).Value = 849.9400000000001
$ws.Cells.Item(5, 6).Value = 0.96
$ws.Cells.Item(5, 7).Value = -43.77

$ws.Cells.Item(6, 2).Value = 150.06
$ws.Cells.Item(6, 5).Value = 1.01
$ws.Cells.Item(6, 6).Value = 0.96
$ws.Cells.Item(6, 7).Value = -8.279999999999999

# --- The remaining per-transaction detail rows (7-11) are no longer     ---
# --- needed once matched against withdrawal dates - drop them, which    ---
# --- shifts the summary rows (---/Gains/Gains excl. losses/Losses) up.  ---
$ws.Range("A7:A11").EntireRow.Delete()

# --- Refresh the summary totals (now rows 8-10) for the new matching ----
$ws.Cells.Item(8, 7).Value = 22.05
$ws.Cells.Item(9, 7).Value = 74.09999999999999
$ws.Cells.Item(10, 7).Value = -52.05

# --- Reflect the updated forex gain/loss total on the ELSTER summary ----
$summary = $wb.Worksheets.Item("ELSTER - Summary")
$summary.Cells.Item(7, 3).Value = 22.05
